# DigitalStewardship.pptx update
#  - give slides 1-6 and 8 an explicit slide background (solid fill ECF0F2)
#  - re-tint slide 7's existing background from E8EEED to ECF0F2
#  - rename three section titles to the "Stewardship" wording and recolor
#    the slide 6 title to match the others

$p = $ppt.ActivePresentation

# RGB(0xEC, 0xF0, 0xF2) packed the way PowerPoint's ColorFormat.RGB expects
# (R + G*256 + B*65536)
$bgColor = 15921388

1..$p.Slides.Count | ForEach-Object {
    $slide = $p.Slides.Item($_)
    $slide.Background.Fill.Solid()
    $slide.Background.Fill.ForeColor.RGB = $bgColor
}

# Slide 4 title: "The Relationship" -> "Stewarding"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Stewarding"

# Slide 5 title: "Imparting Principles" -> "Facilitating Stewardship"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Facilitating Stewardship"

# Slide 6 title: "Conclusion" -> "Encouraging Stewardship", recolored to 26374B
$s6 = $p.Slides.Item(6)
$s6Title = $s6.Shapes.Item(1).TextFrame.TextRange
$s6Title.Text = "Encouraging Stewardship"
$s6Title.Font.Color.RGB = 4929318
